# [BI-1059] Updating files for name length
#
# The "Template" sheet's sample row (row 2) renames the trait:
#   - A2 (Ontology term name): "Powdery Mildew severity field, leaves" -> "PM_Leaf"
#   - J2 (Method description):  "Observed severity of Powdery Mildew on leaves"
#                                -> "Powdery Mildew severity, leaf"
# The active selection/scroll position also moved onto the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$ws.Range("A2").Value = "PM_Leaf"
$ws.Range("J2").Value = "Powdery Mildew severity, leaf"

# Reflect the new view state: scrolled to column E and selection on J2.
$ws.Activate()
$ws.Range("J2").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
